$wb = $excel.ActiveWorkbook

# --- "calciner" sheet (xl/worksheets/sheet9.xml) -----------------------
# Update existing values and append a new "reactor volume" row.
$calciner = $wb.Worksheets.Item("calciner")

$calciner.Range("B3").Value = 200
$calciner.Range("B4").Value = 0.25

$calciner.Range("A5").Value = "reactor volume"
$calciner.Range("B5").Value = 848.2
$calciner.Range("C5").Value = "cm^3"

# Make "calciner" the active/selected sheet (this also clears
# tabSelected on whichever sheet was previously active, i.e.
# "ammonia recovery tower") and move the active cell to B6.
$calciner.Select()
$calciner.Range("B6").Select()
